$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

$ws.Cells.Item(1, 6).Value = "lloq"
$ws.Cells.Item(2, 6).Value = "<0.01"

$ws.Cells.Item(3, 6).NumberFormat = "@"
$ws.Cells.Item(3, 6).Value = "2.23"
$ws.Cells.Item(3, 6).Style = "Normal"

$ws.Range("G11").Select()
